# "first of many passes through in cleaning database"
#
# Column I ("roboticS1Prep") was stored as the text "No" in every data row.
# Convert it to a real boolean (FALSE) value, displayed with a custom
# "TRUE";"TRUE";"FALSE" number format, instead of the shared text string.
# Once no cell references the "No" shared string any more it naturally
# drops out of sharedStrings.xml and every other shared-string index used
# elsewhere on the sheet (rows A14:A49, E8:E49, J2:J49, ...) shifts down to
# stay in sync - no need to touch those cells by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$boolFormat = """TRUE"";""TRUE"";""FALSE"""

for ($r = 2; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 9)   # column I
    $cell.Value = $false
    $cell.NumberFormat = $boolFormat
}

# Mirror the author's final selection: column I (the column just edited)
# instead of column H.
$ws.Range("I2:I49").Select()
